$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H19").Value = 1450.1666
$ws.Range("I19").Value = 1725.25
$ws.Range("J19").Value = 900
$ws.Range("K19").Value = 1725.25
$ws.Range("L19").Value = 900
$ws.Range("M19").Value = -1550.25
$ws.Range("N19").Value = -1250

$ws.Range("H38").Value = 2681.5925
$ws.Range("I38").Value = 88.44444
$ws.Range("J38").Value = 3978.1667
$ws.Range("K38").Value = 265.33332
$ws.Range("L38").Value = 11934.5001
$ws.Range("M38").Value = 106.66668
$ws.Range("N38").Value = -12678.5001

$ws.Range("H42").Value = 160.54546
$ws.Range("I42").Value = 68
$ws.Range("J42").Value = 271.6
$ws.Range("K42").Value = 204
$ws.Range("L42").Value = 814.8000000000001
$ws.Range("M42").Value = 26
$ws.Range("N42").Value = -1274.8

$ws.Range("H137").Value = 2102.1875
$ws.Range("I137").Value = 1604.1818
$ws.Range("K137").Value = 4812.5454
$ws.Range("M137").Value = -2262.5454

$ws.Range("H140").Value = 37550
$ws.Range("J140").Value = 37550
$ws.Range("L140").Value = 37550
$ws.Range("N140").Value = -47910

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H45").Value = 1186
$ws.Range("I45").Value = 1109.25
$ws.Range("J45").Value = 1800
$ws.Range("K45").Value = 1109.25
$ws.Range("L45").Value = 1800
$ws.Range("M45").Value = -732.25
$ws.Range("N45").Value = -2554

$ws.Range("H61").Value = 1650.5652
$ws.Range("I61").Value = 1389
$ws.Range("J61").Value = 2141
$ws.Range("K61").Value = 1389
$ws.Range("L61").Value = 2141
$ws.Range("M61").Value = -1177
$ws.Range("N61").Value = -2565

$ws.Range("H74").Value = 14065315
$ws.Range("I74").Value = 22502868
$ws.Range("J74").Value = 2726.0833
$ws.Range("K74").Value = 22502868
$ws.Range("L74").Value = 2726.0833
$ws.Range("M74").Value = -22501994
$ws.Range("N74").Value = -4474.0833

$ws.Range("H77").Value = 14065315
$ws.Range("I77").Value = 22502868
$ws.Range("J77").Value = 2726.0833
$ws.Range("K77").Value = 112514340
$ws.Range("L77").Value = 13630.4165
$ws.Range("M77").Value = -112509972
$ws.Range("N77").Value = -22366.4165

$ws.Range("H132").Value = 2100.9565
$ws.Range("I132").Value = 1901.4722
$ws.Range("J132").Value = 2819.1
$ws.Range("K132").Value = 5704.4166
$ws.Range("L132").Value = 8457.299999999999
$ws.Range("M132").Value = -3174.4166
$ws.Range("N132").Value = -13517.3

$ws.Range("H136").Value = 1650.5652
$ws.Range("I136").Value = 1389
$ws.Range("J136").Value = 2141
$ws.Range("K136").Value = 4167
$ws.Range("L136").Value = 6423
$ws.Range("M136").Value = -1617
$ws.Range("N136").Value = -11523

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H99").Value = 5000
$ws.Range("I99").Value = 5000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3502
$ws.Range("N99").ClearContents()

$ws.Range("H134").Value = 2185.9062
$ws.Range("I134").Value = 1471.75
$ws.Range("J134").Value = 4328.375
$ws.Range("K134").Value = 4415.25
$ws.Range("L134").Value = 12985.125
$ws.Range("M134").Value = -1880.25
$ws.Range("N134").Value = -18055.125

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H31").Value = 9549097
$ws.Range("I31").Value = 9137467
$ws.Range("J31").Value = 10001891
$ws.Range("K31").Value = 9137467
$ws.Range("L31").Value = 10001891
$ws.Range("N31").Value = -10002481
$ws.Range("M31").Value = -9137172

$ws.Range("H34").Value = 9549097
$ws.Range("I34").Value = 9137467
$ws.Range("J34").Value = 10001891
$ws.Range("K34").Value = 9137467
$ws.Range("L34").Value = 10001891
$ws.Range("N34").Value = -10002295
$ws.Range("M34").Value = -9137265

$ws.Range("H58").Value = 1391.6
$ws.Range("I58").Value = 818.26666
$ws.Range("J58").Value = 1964.9333
$ws.Range("K58").Value = 818.26666
$ws.Range("L58").Value = 1964.9333
$ws.Range("M58").Value = -615.26666
$ws.Range("N58").Value = -2370.9333

$ws.Range("H99").Value = 7153346
$ws.Range("I99").Value = 8938432
$ws.Range("J99").Value = 13000
$ws.Range("K99").Value = 8938432
$ws.Range("L99").Value = 13000
$ws.Range("M99").Value = -8936934
$ws.Range("N99").Value = -15996

$ws.Range("H126").Value = 7153346
$ws.Range("I126").Value = 8938432
$ws.Range("J126").Value = 13000
$ws.Range("K126").Value = 26815296
$ws.Range("L126").Value = 39000
$ws.Range("M126").Value = -26812826
$ws.Range("N126").Value = -43940

$ws.Range("H132").Value = 1736.9
$ws.Range("I132").Value = 1449.6666
$ws.Range("J132").Value = 2598.6
$ws.Range("K132").Value = 4348.9998
$ws.Range("L132").Value = 7795.799999999999
$ws.Range("M132").Value = -1818.9998
$ws.Range("N132").Value = -12855.8

$ws.Range("H134").Value = 7295.1113
$ws.Range("I134").Value = 9123.666999999999
$ws.Range("J134").Value = 3638
$ws.Range("K134").Value = 27371.001
$ws.Range("L134").Value = 10914
$ws.Range("M134").Value = -24836.001
$ws.Range("N134").Value = -15984

$ws.Range("H136").Value = 1391.6
$ws.Range("I136").Value = 818.26666
$ws.Range("J136").Value = 1964.9333
$ws.Range("K136").Value = 2454.79998
$ws.Range("L136").Value = 5894.7999
$ws.Range("M136").Value = 95.20002000000022
$ws.Range("N136").Value = -10994.7999

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H34").Value = 35586.5
$ws.Range("J34").Value = 35586.5
$ws.Range("L34").Value = 35586.5
$ws.Range("N34").Value = -36122.5

$ws.Range("H76").Value = 35586.5
$ws.Range("J76").Value = 35586.5
$ws.Range("L76").Value = 35586.5
$ws.Range("N76").Value = -36216.5

$ws.Range("H79").Value = 35586.5
$ws.Range("J79").Value = 35586.5
$ws.Range("L79").Value = 35586.5
$ws.Range("N79").Value = -37770.5

$ws.Range("H102").Value = 3604
$ws.Range("I102").Value = 3604
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3604
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1982
$ws.Range("N102").ClearContents()

$ws.Range("H132").Value = 2622.0715
$ws.Range("I132").Value = 1464
$ws.Range("J132").Value = 4166.1665
$ws.Range("K132").Value = 4392
$ws.Range("L132").Value = 12498.4995
$ws.Range("M132").Value = -1862
$ws.Range("N132").Value = -17558.4995

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H22").Value = 558.3
$ws.Range("I22").Value = 474.2
$ws.Range("J22").Value = 642.4
$ws.Range("K22").Value = 474.2
$ws.Range("L22").Value = 642.4
$ws.Range("M22").Value = -179.2
$ws.Range("N22").Value = -1232.4

$ws.Range("H27").Value = 558.3
$ws.Range("I27").Value = 474.2
$ws.Range("J27").Value = 642.4
$ws.Range("K27").Value = 474.2
$ws.Range("L27").Value = 642.4
$ws.Range("M27").Value = -367.2
$ws.Range("N27").Value = -856.4

$ws.Range("H46").Value = 2179.5
$ws.Range("I46").Value = 1101
$ws.Range("J46").Value = 2333.5715
$ws.Range("K46").Value = 1101
$ws.Range("L46").Value = 2333.5715
$ws.Range("M46").Value = -913
$ws.Range("N46").Value = -2709.5715

$ws.Range("H55").Value = 369.18182
$ws.Range("I55").Value = 194.83333
$ws.Range("J55").Value = 578.4
$ws.Range("K55").Value = 194.83333
$ws.Range("L55").Value = 578.4
$ws.Range("M55").Value = -21.83332999999999
$ws.Range("N55").Value = -924.4

$ws.Range("H132").Value = 1570309
$ws.Range("I132").Value = 2614927.2
$ws.Range("J132").Value = 3381.625
$ws.Range("K132").Value = 7844781.600000001
$ws.Range("L132").Value = 10144.875
$ws.Range("M132").Value = -7842251.600000001
$ws.Range("N132").Value = -15204.875

$ws.Range("H136").Value = 3576297.2
$ws.Range("I136").Value = 5214422.5
$ws.Range("J136").Value = 2206.2727
$ws.Range("K136").Value = 15643267.5
$ws.Range("L136").Value = 6618.8181
$ws.Range("M136").Value = -15640717.5
$ws.Range("N136").Value = -11718.8181

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H80").Value = 38675.75
$ws.Range("J80").Value = 38675.75
$ws.Range("L80").Value = 38675.75
$ws.Range("N80").Value = -40671.75

$ws.Range("H83").Value = 38675.75
$ws.Range("J83").Value = 38675.75
$ws.Range("L83").Value = 116027.25
$ws.Range("N83").Value = -126011.25

$ws.Range("H113").Value = 15152305
$ws.Range("I113").Value = 20000800
$ws.Range("J113").Value = 756.25
$ws.Range("K113").Value = 60002400
$ws.Range("L113").Value = 2268.75
$ws.Range("M113").Value = -60000230
$ws.Range("N113").Value = -6608.75

$ws.Range("H122").Value = 1000000000
$ws.Range("I122").Value = 1000000000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3000000000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2999997550
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2131.2766
$ws.Range("I132").Value = 1627.7142
$ws.Range("J132").Value = 3600
$ws.Range("K132").Value = 4883.142599999999
$ws.Range("L132").Value = 10800
$ws.Range("M132").Value = -2353.142599999999
$ws.Range("N132").Value = -15860

$ws.Range("H136").Value = 1824.84
$ws.Range("I136").Value = 1399.4
$ws.Range("J136").Value = 2463
$ws.Range("K136").Value = 4198.200000000001
$ws.Range("L136").Value = 7389
$ws.Range("M136").Value = -1648.200000000001
$ws.Range("N136").Value = -12489
